$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append 9 new master-data rows (rows 22-30), continuing the existing
# regcntr_id / machine_id pattern:
#   A = regcntr_id, B = machine_id, C = lang_code ("eng"),
#   D = is_active (TRUE), E = cr_by ("superadmin"), F = cr_dtimes ("now()")
$newRows = @(
    @(10002, 10021),
    @(10003, 10022),
    @(10004, 10023),
    @(10005, 10024),
    @(10006, 10025),
    @(10007, 10026),
    @(10008, 10027),
    @(10009, 10028),
    @(10010, 10029)
)

$row = 22
foreach ($pair in $newRows) {
    $ws.Cells.Item($row, 1).Value = $pair[0]
    $ws.Cells.Item($row, 2).Value = $pair[1]
    $ws.Cells.Item($row, 3).Value = "eng"
    $ws.Cells.Item($row, 4).Value = $true
    $ws.Cells.Item($row, 5).Value = "superadmin"
    $ws.Cells.Item($row, 6).Value = "now()"
    $row++
}

# Basic print setup tweak that was part of the authored edit.
$ws.PageSetup.Orientation = 1

# Leave the sheet scrolled/selected the way the author left it: viewing the
# newly-added tail of the table with the row right after the data selected
# (selecting the remainder of the sheet below the table, like pressing
# Ctrl+Shift+End-style "select to bottom" after the last row).
$ws.Range("A31:XFD1048576").Select()
